# Age distribution. Fix #51.
# The "help_age_content" string (shown as help text under the age question)
# is expanded to mention the age range (57-95 years) covered by the
# underlying study data, in both the Swedish (column B) and English
# (column C) localisations, on row 16 of Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B16").Value = "Hur gammal kommer du att vara vid den planerade operationen? <p> I vårt underliggande studiematerial var den yngste patienten som dog 57 år och den äldste 95 år. Vi har således inga tillförlitliga data för skattningar utanför detta intervall."
$ws.Range("C16").Value = "How old will you be at the time of your planned hip replacement? <p> The youngest and oldest patients who died in the study cohort were 57 and 95 years old. We are unable to provide estimtes outside this range."

# The longer help text now wraps onto two lines, so the row grows taller
# (matching the other two-line help rows such as 18 and 26).
$ws.Rows.Item(16).RowHeight = 28.8

# The active selection in the saved file moved to C16.
$ws.Range("C16").Select()
